$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -71.3861
$ws.Range("B2").Value = -71.3798
$ws.Range("A3").Value = 41.562
$ws.Range("B3").Value = 41.5666
$ws.Range("A4").Value = -71.2835
$ws.Range("B4").Value = -71.2898
$ws.Range("A5").Value = 41.6852
$ws.Range("B5").Value = 41.6805
